$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A93").Value = "overworld_1_latitude_0"
$ws.Range("A94").Value = "overworld_1_temp_intro_0"

$ws.Range("B82").Value = "It looks like they have been exiled from their planet and are looking for a new home."
$ws.Range("B84").Value = "They are expressing their gratitude and are ready to cooperate."
$ws.Range("B105").Value = "Since the earth rotates at a slightly tilted axis around the Sun, the atmosphere can change throughout the year."
$ws.Range("B106").Value = "Let's adjust the time by selecting a different season. Perhaps summer will give us the temperature to satisfy the criteria!"
$ws.Range("B132").Value = "Though the weather is ideal for the frogs, unfortunately, it is also ideal for these invasive plants."
$ws.Range("B152").Value = "This time around, there is more than one hotspot to discover on the map. Only one of them will match with the frogs’ criteria."
$ws.Range("B153").Value = "Looks like we’ve landed in a tropical climate, where it’s hot and humid all year round with plenty of rain."
$ws.Range("B159").Value = "Look out! It's a beetle!"
$ws.Range("B161").Value = "These troublesome insects can be dealt with by a hero frog. Make sure to have one around to get them out."
$ws.Range("B162").Value = "Take cover! A hurricane is heading our way!"

$ws.Activate()
$ws.Range("A94").Select()
$excel.ActiveWindow.ScrollRow = 84
